$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain string/text updates (not numeric-looking, so Excel keeps them as text)
$ws.Range("D2").Value = "62.014.13"
$ws.Range("E2").Value = "  -4.45%  "
$ws.Range("D3").Value = "2.999.27"
$ws.Range("E3").Value = "  -5.60%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("E5").Value = "  -3.88%  "
$ws.Range("E6").Value = "  -6.94%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").Value = "  -4.57%  "
$ws.Range("D9").Value = "3.006.71"
$ws.Range("E9").Value = "  -5.11%  "
$ws.Range("E10").Value = "  -5.06%  "
$ws.Range("E11").Value = "  -6.61%  "
$ws.Range("E12").Value = "  -4.88%  "
$ws.Range("D13").Value = "3.525.24"
$ws.Range("E13").Value = "  -5.32%  "
$ws.Range("E14").Value = "  -3.78%  "
$ws.Range("D15").Value = "62.098.03"
$ws.Range("E15").Value = "  -4.20%  "
$ws.Range("E16").Value = "  -6.88%  "
$ws.Range("D17").Value = "3.000.57"
$ws.Range("E17").Value = "  -5.76%  "
$ws.Range("E18").Value = "  -5.08%  "
$ws.Range("E19").Value = "  -4.22%  "
$ws.Range("E20").Value = "  -3.03%  "
$ws.Range("E21").Value = "  -5.95%  "
$ws.Range("E22").Value = "  -6.74%  "
$ws.Range("E23").Value = "  -0.14%  "
$ws.Range("E24").Value = "  -5.05%  "
$ws.Range("E25").Value = "  -3.73%  "
$ws.Range("E26").Value = "  -7.71%  "
$ws.Range("E27").Value = "  -9.02%  "
$ws.Range("E28").Value = "  +0.05%  "
$ws.Range("E29").Value = "  -4.47%  "
$ws.Range("E30").Value = "  -0.03%  "
$ws.Range("E31").Value = "  -5.54%  "
$ws.Range("E32").Value = "  -4.19%  "
$ws.Range("E33").Value = "  +2.40%  "
$ws.Range("E34").Value = "  -6.08%  "
$ws.Range("E35").Value = "  -4.81%  "
$ws.Range("E36").Value = "  -4.67%  "
$ws.Range("E37").Value = "  -5.30%  "
$ws.Range("E38").Value = "  -8.95%  "
$ws.Range("D39").Value = "2.452.09"
$ws.Range("E39").Value = "  -10.32%  "
$ws.Range("E40").Value = "  -4.25%  "
$ws.Range("E41").Value = "  -6.07%  "
$ws.Range("E42").Value = "  -4.39%  "
$ws.Range("E43").Value = "  -6.00%  "
$ws.Range("E44").Value = "  -6.22%  "
$ws.Range("E45").Value = "  -0.16%  "
$ws.Range("E46").Value = "  -5.48%  "
$ws.Range("E47").Value = "  -10.72%  "
$ws.Range("E48").Value = "  -6.97%  "
$ws.Range("E49").Value = "  -3.76%  "
$ws.Range("E50").Value = "  +0.32%  "
$ws.Range("E51").Value = "  -9.08%  "

# Numeric-looking price strings: write as a formula producing the literal text,
# then paste-special as values so the cell ends up a plain text cell (matches the
# original inlineStr cells) instead of Excel auto-converting the text to a number.
$ws.Range("D5").Formula = "=""550.42"""
$ws.Range("D5").Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("D6").Formula = "=""154.67"""
$ws.Range("D6").Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("D11").Formula = "=""6.25"""
$ws.Range("D11").Copy()
$ws.Range("D11").PasteSpecial(-4163)
$ws.Range("D19").Formula = "=""393.59"""
$ws.Range("D19").Copy()
$ws.Range("D19").PasteSpecial(-4163)
$ws.Range("D21").Formula = "=""11.96"""
$ws.Range("D21").Copy()
$ws.Range("D21").PasteSpecial(-4163)
$ws.Range("D22").Formula = "=""6.65"""
$ws.Range("D22").Copy()
$ws.Range("D22").PasteSpecial(-4163)
$ws.Range("D24").Formula = "=""65.17"""
$ws.Range("D24").Copy()
$ws.Range("D24").PasteSpecial(-4163)
$ws.Range("D25").Formula = "=""0.468"""
$ws.Range("D25").Copy()
$ws.Range("D25").PasteSpecial(-4163)
$ws.Range("D26").Formula = "=""0.185"""
$ws.Range("D26").Copy()
$ws.Range("D26").PasteSpecial(-4163)
$ws.Range("D28").Formula = "=""0.998"""
$ws.Range("D28").Copy()
$ws.Range("D28").PasteSpecial(-4163)
$ws.Range("D29").Formula = "=""8.46"""
$ws.Range("D29").Copy()
$ws.Range("D29").PasteSpecial(-4163)
$ws.Range("D32").Formula = "=""20.48"""
$ws.Range("D32").Copy()
$ws.Range("D32").PasteSpecial(-4163)
$ws.Range("D33").Formula = "=""159.37"""
$ws.Range("D33").Copy()
$ws.Range("D33").PasteSpecial(-4163)
$ws.Range("D40").Formula = "=""3.92"""
$ws.Range("D40").Copy()
$ws.Range("D40").PasteSpecial(-4163)
$ws.Range("D41").Formula = "=""22.41"""
$ws.Range("D41").Copy()
$ws.Range("D41").PasteSpecial(-4163)
$ws.Range("D42").Formula = "=""37.23"""
$ws.Range("D42").Copy()
$ws.Range("D42").PasteSpecial(-4163)
$ws.Range("D45").Formula = "=""0.998"""
$ws.Range("D45").Copy()
$ws.Range("D45").PasteSpecial(-4163)
$ws.Range("D47").Formula = "=""4.97"""
$ws.Range("D47").Copy()
$ws.Range("D47").PasteSpecial(-4163)
$ws.Range("D48").Formula = "=""19.86"""
$ws.Range("D48").Copy()
$ws.Range("D48").PasteSpecial(-4163)
$ws.Range("D51").Formula = "=""265.97"""
$ws.Range("D51").Copy()
$ws.Range("D51").PasteSpecial(-4163)

$excel.CutCopyMode = $false
